$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8,9,10 (ECs sending-cluster rows disappear from the
# dataset; overall row count drops from 10 to 7)
$ws.Rows.Item(8).Resize(3).Delete() | Out-Null

# New data (updated TPM-derived values) for rows 2-7.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("FAPs","Alcam","Nrp1","ECs",2,0.6666666666666666,0.4771596666666666,1.431479,0.4973652976730675,0.4973652976730676,3,1,123.2806423333333,369.841927,0.6241574062367528,0.6241574062367526,58.82455020222589,529.420951820033,0.3104342341477923,0.3104342341477923),
    @("FAPs","Alcam","Nrp1","FAPs",2,0.6666666666666666,0.4771596666666666,1.431479,0.4973652976730675,0.4973652976730676,3,1,47.26005833333333,141.780175,0.2392728888301323,0.2392728888301322,22.55059368098055,202.955343128825,0.1190060315780935,0.1190060315780935),
    @("FAPs","Alcam","Nrp1","MuSCs",2,0.6666666666666666,0.4771596666666666,1.431479,0.4973652976730675,0.4973652976730676,3,1,26.97460733333333,80.923822,0.136569704933115,0.136569704933115,12.87119464363756,115.840751792738,0.06792503194718176,0.06792503194718176),
    @("MuSCs","Alcam","Nrp1","ECs",3,1,0.482215,1.446645,0.5026347023269324,0.5026347023269324,3,1,123.2806423333333,369.841927,0.6241574062367528,0.6241574062367526,59.44777494276833,535.029974484915,0.3137231720889604,0.3137231720889604),
    @("MuSCs","Alcam","Nrp1","FAPs",3,1,0.482215,1.446645,0.5026347023269324,0.5026347023269324,3,1,47.26005833333333,141.780175,0.2392728888301323,0.2392728888301322,22.78950902920833,205.105581262875,0.1202668572520387,0.1202668572520387),
    @("MuSCs","Alcam","Nrp1","MuSCs",3,1,0.482215,1.446645,0.5026347023269324,0.5026347023269324,3,1,26.97460733333333,80.923822,0.136569704933115,0.136569704933115,13.00756027524333,117.06804247719,0.06864467298593327,0.06864467298593326)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Output "done"
